$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E range to Text format before writing so numeric-looking
# strings (e.g. "0.9983") are preserved as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.151.87'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '1.767.01'
$ws.Range("E3").Value = '  +2.96%  '

$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("D5").Value = '312.73'
$ws.Range("E5").Value = '  +1.38%  '

$ws.Range("D6").Value = '0.9979'
$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("D7").Value = '0.5215'
$ws.Range("E7").Value = '  +10.60%  '

$ws.Range("D8").Value = '0.3604'
$ws.Range("E8").Value = '  +5.14%  '

$ws.Range("D9").Value = '42.32'
$ws.Range("E9").Value = '  +0.45%  '

$ws.Range("D10").Value = '0.07327'
$ws.Range("E10").Value = '  +0.80%  '

$ws.Range("D11").Value = '1.079'
$ws.Range("E11").Value = '  +3.49%  '

$ws.Range("D12").Value = '0.9979'
$ws.Range("E12").Value = '  -0.37%  '

$ws.Range("E13").Value = '  +3.20%  '

$ws.Range("D14").Value = '6.047'
$ws.Range("E14").Value = '  +3.05%  '

$ws.Range("D15").Value = '1.765.53'
$ws.Range("E15").Value = '  +2.44%  '

$ws.Range("D16").Value = '6.945'
$ws.Range("E16").Value = '  +0.90%  '

$ws.Range("D17").Value = '88.20'
$ws.Range("E17").Value = '  -1.13%  '

$ws.Range("D18").Value = '0.00001041'
$ws.Range("E18").Value = '  +0.16%  '

$ws.Range("D19").Value = '0.06412'

$ws.Range("D20").Value = '0.9977'
$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("D21").Value = '16.71'
$ws.Range("E21").Value = '  +1.20%  '

$ws.Range("D22").Value = '5.823'
$ws.Range("E22").Value = '  +3.57%  '

$ws.Range("D23").Value = '27.246.45'
$ws.Range("E23").Value = '  +0.41%  '

$ws.Range("D24").Value = '11.34'
$ws.Range("E24").Value = '  +4.28%  '

$ws.Range("D25").Value = '2.064'
$ws.Range("E25").Value = '  -2.71%  '

$ws.Range("D26").Value = '154.35'
$ws.Range("E26").Value = '  -1.17%  '

$ws.Range("D27").Value = '20.05'
$ws.Range("E27").Value = '  +2.78%  '

$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '1.963.27'
$ws.Range("E28").Value = '  +2.72%  '

$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '2.304'
$ws.Range("E29").Value = '  +10.28%  '

$ws.Range("D30").Value = '120.83'
$ws.Range("E30").Value = '  +1.08%  '

$ws.Range("D31").Value = '1.056'
$ws.Range("E31").Value = '  +4.04%  '

$ws.Range("D32").Value = '0.09722'
$ws.Range("E32").Value = '  +6.24%  '

$ws.Range("D33").Value = '5.494'
$ws.Range("E33").Value = '  +3.22%  '

$ws.Range("D34").Value = '3.595'
$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").Value = '0.02216'
$ws.Range("E35").Value = '  +0.60%  '

$ws.Range("D36").Value = '0.05977'

$ws.Range("D37").Value = '11.18'
$ws.Range("E37").Value = '  +1.58%  '

$ws.Range("D38").Value = '0.2026'
$ws.Range("E38").Value = '  +1.45%  '

$ws.Range("D39").Value = '4.820'
$ws.Range("E39").Value = '  +1.78%  '

$ws.Range("D40").Value = '0.6112'
$ws.Range("E40").Value = '  +3.75%  '

$ws.Range("D41").Value = '1.432'
$ws.Range("E41").Value = '  +2.66%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '7.922'
$ws.Range("E42").Value = '  +6.03%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.125'
$ws.Range("E43").Value = '  +0.73%  '

$ws.Range("D44").Value = '13.19'
$ws.Range("E44").Value = '  +5.10%  '

$ws.Range("D45").Value = '3.618'
$ws.Range("E45").Value = '  +1.52%  '

$ws.Range("D46").Value = '0.5731'
$ws.Range("E46").Value = '  +1.41%  '

$ws.Range("D47").Value = '121.17'
$ws.Range("E47").Value = '  +3.25%  '

$ws.Range("D48").Value = '1.875'
$ws.Range("E48").Value = '  +1.80%  '

$ws.Range("D49").Value = '1.107'
$ws.Range("E49").Value = '  +1.86%  '

$ws.Range("D50").Value = '0.06694'
$ws.Range("E50").Value = '  +0.65%  '

$ws.Range("D51").Value = '70.69'
$ws.Range("E51").Value = '  +1.30%  '

# Restore the default "Normal" style so no stray per-cell number format
# is left behind (keeps the written cells text-typed, style unchanged).
$ws.Range("D2:E51").Style = "Normal"
